$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 5
    3  = 2
    4  = 2
    5  = 3
    6  = 5
    7  = 9
    8  = 8
    9  = 4
    10 = 3
    11 = 4
    12 = 7
    13 = 6
    14 = 8
    15 = 10
    16 = 10
    17 = 8
    18 = 7
    19 = 9
    20 = 8
    21 = 7
    22 = 9
    23 = 7
    24 = 11
    25 = 5
    26 = 11
    27 = 7
    28 = 9
    29 = 8
    30 = 7
    31 = 6
    32 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}

$wb.Save()
